$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.177.47'
$ws.Range("E2").Value = '  +0.19%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.423.41'
$ws.Range("E3").Value = '  +0.14%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '554.38'
$ws.Range("E5").Value = '  +0.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.23'
$ws.Range("E6").Value = '  -0.34%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.589'
$ws.Range("E8").Value = '  +2.14%  '
$ws.Range("E9").Value = '  -1.23%  '
$ws.Range("E10").Value = '  -0.78%  '
$ws.Range("E11").Value = '  -0.26%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.354'
$ws.Range("E12").Value = '  -1.48%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.92'
$ws.Range("E13").Value = '  +0.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.852.91'
$ws.Range("E14").Value = '  +0.11%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.090.54'
$ws.Range("E15").Value = '  +0.19%  '
$ws.Range("E16").Value = '  -0.21%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.422.22'
$ws.Range("E17").Value = '  +0.57%  '
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("E19").Value = '  +2.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '327.25'
$ws.Range("E20").Value = '  -1.57%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.76'
$ws.Range("E21").Value = '  +0.25%  '
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.26'
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("E24").Value = '  +4.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.66'
$ws.Range("E25").Value = '  +0.86%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.42'
$ws.Range("E27").Value = '  +5.73%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0775'
$ws.Range("E28").Value = '  -0.97%  '
$ws.Range("E29").Value = '  -0.23%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '170.63'
$ws.Range("E30").Value = '  +0.86%  '
$ws.Range("E31").Value = '  -2.38%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.403'
$ws.Range("E32").Value = '  -3.62%  '
$ws.Range("E33").Value = '  +2.30%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.57'
$ws.Range("E34").Value = '  -0.57%  '
$ws.Range("E35").Value = '  +3.12%  '
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.23'
$ws.Range("E37").Value = '  +0.43%  '
$ws.Range("E38").Value = '  +0.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '327.15'
$ws.Range("E39").Value = '  +4.24%  '
$ws.Range("E40").Value = '  -0.48%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '145.34'
$ws.Range("E41").Value = '  +4.34%  '
$ws.Range("E42").Value = '  -0.83%  '
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '19.93'
$ws.Range("E43").Value = '  +2.11%  '
$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0964'
$ws.Range("E44").Value = '  +0.36%  '
$ws.Range("E45").Value = '  -0.70%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.576'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '11.04'
$ws.Range("E48").Value = '  -0.06%  '
$ws.Range("E49").Value = '  -1.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.67'
$ws.Range("E50").Value = '  -0.48%  '
$ws.Range("E51").Value = '  -0.65%  '
